# Auto-generated edit script applying the cryptos.xlsx price/volume/hour refresh
# described by the commit "Updated symbol list on Fri Jan 27 14:12:44 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is the target cell and its new literal text value. Values such as
# "304.49", "-0.82%" or "14" look numeric/percentage to Excel, so without forcing
# the cell to Text format first, Excel would silently coerce them into numbers
# (e.g. 304.49 -> 304.49000000000001, "-0.82%" -> -0.0082) instead of keeping the
# plain text the source workbook stores them as.
$updates = @(
    @{ Cell = 'D2'; Value = '304.49' }
    @{ Cell = 'E2'; Value = '-0.82%' }
    @{ Cell = 'G2'; Value = '14' }
    @{ Cell = 'D3'; Value = '35.75' }
    @{ Cell = 'E3'; Value = '-0.32%' }
    @{ Cell = 'G3'; Value = '14' }
    @{ Cell = 'D4'; Value = '5.034' }
    @{ Cell = 'E4'; Value = '-0.64%' }
    @{ Cell = 'G4'; Value = '14' }
    @{ Cell = 'D5'; Value = '0.07971' }
    @{ Cell = 'E5'; Value = '-1.76%' }
    @{ Cell = 'G5'; Value = '14' }
    @{ Cell = 'D6'; Value = '1.860' }
    @{ Cell = 'E6'; Value = '-4.21%' }
    @{ Cell = 'G6'; Value = '14' }
    @{ Cell = 'B7'; Value = 'GateToken' }
    @{ Cell = 'C7'; Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt' }
    @{ Cell = 'D7'; Value = '4.120' }
    @{ Cell = 'E7'; Value = '-0.42%' }
    @{ Cell = 'G7'; Value = '14' }
    @{ Cell = 'B8'; Value = 'KuCoinToken' }
    @{ Cell = 'C8'; Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs' }
    @{ Cell = 'D8'; Value = '7.775' }
    @{ Cell = 'E8'; Value = '-0.03%' }
    @{ Cell = 'G8'; Value = '14' }
    @{ Cell = 'B9'; Value = 'MXToken' }
    @{ Cell = 'C9'; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx' }
    @{ Cell = 'D9'; Value = '0.9209' }
    @{ Cell = 'E9'; Value = '-1.12%' }
    @{ Cell = 'G9'; Value = '14' }
    @{ Cell = 'B10'; Value = 'LiechtensteinCryptoassetsExchange' }
    @{ Cell = 'C10'; Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx' }
    @{ Cell = 'D10'; Value = '0.1277' }
    @{ Cell = 'E10'; Value = '-4.93%' }
    @{ Cell = 'G10'; Value = '14' }
    @{ Cell = 'B11'; Value = 'WazirX' }
    @{ Cell = 'C11'; Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx' }
    @{ Cell = 'D11'; Value = '0.1887' }
    @{ Cell = 'E11'; Value = '-1.52%' }
    @{ Cell = 'G11'; Value = '14' }
    @{ Cell = 'B12'; Value = 'MandalaExchangeToken' }
    @{ Cell = 'C12'; Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx' }
    @{ Cell = 'D12'; Value = '0.09053' }
    @{ Cell = 'E12'; Value = '-1.95%' }
    @{ Cell = 'G12'; Value = '14' }
    @{ Cell = 'B13'; Value = 'BitrueCoin' }
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr' }
    @{ Cell = 'D13'; Value = '0.03420' }
    @{ Cell = 'E13'; Value = '-2.17%' }
    @{ Cell = 'G13'; Value = '14' }
    @{ Cell = 'B14'; Value = 'BitMartToken' }
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx' }
    @{ Cell = 'D14'; Value = '0.09850' }
    @{ Cell = 'E14'; Value = '-0.28%' }
    @{ Cell = 'G14'; Value = '14' }
    @{ Cell = 'B15'; Value = 'BitForexToken' }
    @{ Cell = 'C15'; Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf' }
    @{ Cell = 'D15'; Value = '0.001403' }
    @{ Cell = 'E15'; Value = '-0.65%' }
    @{ Cell = 'G15'; Value = '14' }
    @{ Cell = 'B16'; Value = 'CoinExToken' }
    @{ Cell = 'C16'; Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet' }
    @{ Cell = 'D16'; Value = '0.04422' }
    @{ Cell = 'E16'; Value = '0.52%' }
    @{ Cell = 'G16'; Value = '14' }
    @{ Cell = 'B17'; Value = 'TigerCash' }
    @{ Cell = 'C17'; Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch' }
    @{ Cell = 'D17'; Value = '0.006202' }
    @{ Cell = 'E17'; Value = '6.81%' }
    @{ Cell = 'G17'; Value = '14' }
    @{ Cell = 'B18'; Value = 'LEO' }
    @{ Cell = 'C18'; Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo' }
    @{ Cell = 'D18'; Value = '3.850' }
    @{ Cell = 'E18'; Value = '6.88%' }
    @{ Cell = 'G18'; Value = '14' }
    @{ Cell = 'B19'; Value = 'BTSEToken' }
    @{ Cell = 'C19'; Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse' }
    @{ Cell = 'D19'; Value = '3.320' }
    @{ Cell = 'E19'; Value = '13.40%' }
    @{ Cell = 'G19'; Value = '14' }
    @{ Cell = 'B20'; Value = 'BitpandaEcosystemToken' }
    @{ Cell = 'C20'; Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best' }
    @{ Cell = 'D20'; Value = '0.3407' }
    @{ Cell = 'E20'; Value = '-1.11%' }
    @{ Cell = 'G20'; Value = '14' }
    @{ Cell = 'B21'; Value = 'ProBitToken' }
    @{ Cell = 'C21'; Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob' }
    @{ Cell = 'D21'; Value = '0.1340' }
    @{ Cell = 'E21'; Value = '0.64%' }
    @{ Cell = 'G21'; Value = '14' }
    @{ Cell = 'B22'; Value = 'MCDex' }
    @{ Cell = 'C22'; Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb' }
    @{ Cell = 'D22'; Value = '4.796' }
    @{ Cell = 'E22'; Value = '-8.40%' }
    @{ Cell = 'G22'; Value = '14' }
    @{ Cell = 'B23'; Value = 'ZBToken' }
    @{ Cell = 'C23'; Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb' }
    @{ Cell = 'D23'; Value = '0.2499' }
    @{ Cell = 'E23'; Value = '-3.65%' }
    @{ Cell = 'G23'; Value = '14' }
    @{ Cell = 'D24'; Value = '0.001233' }
    @{ Cell = 'E24'; Value = '1.08%' }
    @{ Cell = 'G24'; Value = '14' }
    @{ Cell = 'E25'; Value = '1.66%' }
    @{ Cell = 'G25'; Value = '14' }
    @{ Cell = 'G26'; Value = '14' }
    @{ Cell = 'E27'; Value = '-21.21%' }
    @{ Cell = 'G27'; Value = '14' }
    @{ Cell = 'E28'; Value = '42.46%' }
    @{ Cell = 'G28'; Value = '14' }
    @{ Cell = 'G29'; Value = '14' }
    @{ Cell = 'G30'; Value = '14' }
    @{ Cell = 'G31'; Value = '14' }
    @{ Cell = 'G32'; Value = '14' }
    @{ Cell = 'G33'; Value = '14' }
    @{ Cell = 'G34'; Value = '14' }
    @{ Cell = 'G35'; Value = '14' }
    @{ Cell = 'G36'; Value = '14' }
    @{ Cell = 'G37'; Value = '14' }
    @{ Cell = 'G38'; Value = '14' }
    @{ Cell = 'D39'; Value = '0.01924' }
    @{ Cell = 'E39'; Value = '-3.58%' }
    @{ Cell = 'G39'; Value = '14' }
    @{ Cell = 'D40'; Value = '0.05145' }
    @{ Cell = 'E40'; Value = '1.77%' }
    @{ Cell = 'G40'; Value = '14' }
    @{ Cell = 'D41'; Value = '0.007551' }
    @{ Cell = 'E41'; Value = '-0.87%' }
    @{ Cell = 'G41'; Value = '14' }
    @{ Cell = 'D42'; Value = '0.01013' }
    @{ Cell = 'E42'; Value = '-9.67%' }
    @{ Cell = 'G42'; Value = '14' }
    @{ Cell = 'D43'; Value = '0.1346' }
    @{ Cell = 'E43'; Value = '-2.45%' }
    @{ Cell = 'G43'; Value = '14' }
    @{ Cell = 'E44'; Value = '0.75%' }
    @{ Cell = 'G44'; Value = '14' }
    @{ Cell = 'D45'; Value = '0.009864' }
    @{ Cell = 'E45'; Value = '-12.75%' }
    @{ Cell = 'G45'; Value = '14' }
    @{ Cell = 'D46'; Value = '0.00006193' }
    @{ Cell = 'E46'; Value = '-3.25%' }
    @{ Cell = 'G46'; Value = '14' }
    @{ Cell = 'E47'; Value = '0.35%' }
    @{ Cell = 'G47'; Value = '14' }
    @{ Cell = 'D48'; Value = '64.99' }
    @{ Cell = 'E48'; Value = '-0.35%' }
    @{ Cell = 'G48'; Value = '14' }
    @{ Cell = 'D49'; Value = '0.001251' }
    @{ Cell = 'E49'; Value = '5.33%' }
    @{ Cell = 'G49'; Value = '14' }
    @{ Cell = 'E50'; Value = '0.35%' }
    @{ Cell = 'G50'; Value = '14' }
    @{ Cell = 'E51'; Value = '0.35%' }
    @{ Cell = 'G51'; Value = '14' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.Style = "Normal"
}

